$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 286010

$ws.Range("B15").Value = 152
$ws.Range("C15").Value = "Cluj-Cmp. Turzii"

$ws.Range("B16").Value = 47
$ws.Range("C16").Value = "Cluj-Cluj"
$ws.Range("D16").Value = "Interes Serviciu"

$ws.Range("B19").Value = 121
$ws.Range("C19").Value = "Cluj-Turda"
$ws.Range("D19").Value = "Interes Serviciu"

$ws.Range("B20").Value = 30
$ws.Range("C20").Value = "Acasa-Birou"
$ws.Range("D20").Value = " "

$ws.Range("B21").Value = 92
$ws.Range("C21").Value = "Cluj-Bontida"

$ws.Range("B22").Value = 121
$ws.Range("C22").Value = "Cluj-Turda"
$ws.Range("D22").Value = "Interes Serviciu"

$ws.Range("B27").Value = 47
$ws.Range("C27").Value = "Cluj-Cluj"
$ws.Range("D27").Value = "Interes Serviciu"

$ws.Range("B28").Value = 30
$ws.Range("C28").Value = "Acasa-Birou"
$ws.Range("D28").Value = " "

$ws.Range("B29").Value = 356
$ws.Range("C29").Value = "Cluj-Baia-Mare"
$ws.Range("D29").Value = "Interes Serviciu"

$ws.Range("B30").Value = 30
$ws.Range("C30").Value = "Acasa-Birou"
$ws.Range("D30").Value = " "

$ws.Range("B33").Value = 121
$ws.Range("C33").Value = "Cluj-Turda"

$ws.Range("B34").Value = 257
$ws.Range("C34").Value = "Cluj-Bistrita"
$ws.Range("D34").Value = "Interes Serviciu"

$ws.Range("B35").Value = 152
$ws.Range("C35").Value = "Cluj-Cmp. Turzii"
$ws.Range("D35").Value = "Interes Serviciu"

$ws.Range("B36").Value = 85
$ws.Range("C36").Value = "Cluj-Apahida"
$ws.Range("D36").Value = "Interes Serviciu"

$ws.Range("B37").Value = 92
$ws.Range("C37").Value = "Cluj-Bontida"

$ws.Range("B40").Value = 121
$ws.Range("C40").Value = "Cluj-Turda"

$ws.Range("B41").Value = 156
$ws.Range("C41").Value = "Cluj-Zalau"
$ws.Range("D41").Value = "Interes Serviciu"

$ws.Range("B42").Value = 152
$ws.Range("C42").Value = "Cluj-Cmp. Turzii"
$ws.Range("D42").Value = "Interes Serviciu"

$ws.Range("B44").Value = 2222
$ws.Range("B45").Value = 288232
